# Update cryptos list values (price & volume columns) to reflect the
# scheduled GitHub Actions refresh. Also renames row 51 from Stellar to SuiNetwork.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.111.81"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "3.517.00"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +4.32%  "
$ws.Range("E9").Value = "  +6.58%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.437"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "4.124.71"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "67.112.65"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "3.509.91"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "394.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  +4.20%  "
$ws.Range("D44").Value = "2.804.58"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "336.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.849"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
